# Auto-generated update of cryptocurrency price/volume/hour data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings
# (prices, percentages, hour) are preserved as text, matching the
# original inline-string cell type rather than being coerced to numbers.
$dataRange = $ws.Range("D2:G51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "301.23"
$ws.Range("E2").Value = "-0.79%"
$ws.Range("G2").Value = "21"
$ws.Range("D3").Value = "31.38"
$ws.Range("E3").Value = "-2.27%"
$ws.Range("G3").Value = "21"
$ws.Range("D4").Value = "5.108"
$ws.Range("E4").Value = "-3.11%"
$ws.Range("G4").Value = "21"
$ws.Range("D5").Value = "0.07366"
$ws.Range("E5").Value = "-1.93%"
$ws.Range("G5").Value = "21"
$ws.Range("D6").Value = "2.399"
$ws.Range("E6").Value = "57.26%"
$ws.Range("G6").Value = "21"
$ws.Range("D7").Value = "7.950"
$ws.Range("E7").Value = "0.96%"
$ws.Range("G7").Value = "21"
$ws.Range("D8").Value = "3.795"
$ws.Range("E8").Value = "-0.66%"
$ws.Range("G8").Value = "21"
$ws.Range("D9").Value = "0.9177"
$ws.Range("E9").Value = "-0.34%"
$ws.Range("G9").Value = "21"
$ws.Range("D10").Value = "0.1714"
$ws.Range("E10").Value = "0.86%"
$ws.Range("G10").Value = "21"
$ws.Range("D11").Value = "0.07632"
$ws.Range("E11").Value = "-3.29%"
$ws.Range("G11").Value = "21"
$ws.Range("D12").Value = "0.08102"
$ws.Range("E12").Value = "1.83%"
$ws.Range("G12").Value = "21"
$ws.Range("D13").Value = "0.03020"
$ws.Range("E13").Value = "-1.85%"
$ws.Range("G13").Value = "21"
$ws.Range("D14").Value = "0.09933"
$ws.Range("E14").Value = "0.37%"
$ws.Range("G14").Value = "21"
$ws.Range("D15").Value = "0.001494"
$ws.Range("E15").Value = "0.11%"
$ws.Range("G15").Value = "21"
$ws.Range("D16").Value = "0.006159"
$ws.Range("E16").Value = "-1.96%"
$ws.Range("G16").Value = "21"
$ws.Range("D17").Value = "3.454"
$ws.Range("E17").Value = "-0.15%"
$ws.Range("G17").Value = "21"
$ws.Range("D18").Value = "2.226"
$ws.Range("E18").Value = "-0.23%"
$ws.Range("G18").Value = "21"
$ws.Range("D19").Value = "0.3295"
$ws.Range("E19").Value = "-0.26%"
$ws.Range("G19").Value = "21"
$ws.Range("E20").Value = "-0.54%"
$ws.Range("G20").Value = "21"
$ws.Range("D21").Value = "4.648"
$ws.Range("E21").Value = "3.18%"
$ws.Range("G21").Value = "21"
$ws.Range("D22").Value = "0.04643"
$ws.Range("E22").Value = "0.77%"
$ws.Range("G22").Value = "21"
$ws.Range("E23").Value = "-3.27%"
$ws.Range("G23").Value = "21"
$ws.Range("E24").Value = "0.81%"
$ws.Range("G24").Value = "21"
$ws.Range("D25").Value = "0.004485"
$ws.Range("E25").Value = "0.91%"
$ws.Range("G25").Value = "21"
$ws.Range("D26").Value = "0.0001298"
$ws.Range("E26").Value = "-6.97%"
$ws.Range("G26").Value = "21"
$ws.Range("D27").Value = "0.0002667"
$ws.Range("E27").Value = "49.87%"
$ws.Range("G27").Value = "21"
$ws.Range("G28").Value = "21"
$ws.Range("G29").Value = "21"
$ws.Range("G30").Value = "21"
$ws.Range("G31").Value = "21"
$ws.Range("G32").Value = "21"
$ws.Range("G33").Value = "21"
$ws.Range("G34").Value = "21"
$ws.Range("G35").Value = "21"
$ws.Range("G36").Value = "21"
$ws.Range("G37").Value = "21"
$ws.Range("G38").Value = "21"
$ws.Range("D39").Value = "0.01738"
$ws.Range("E39").Value = "1.16%"
$ws.Range("G39").Value = "21"
$ws.Range("D40").Value = "0.04524"
$ws.Range("E40").Value = "0.90%"
$ws.Range("G40").Value = "21"
$ws.Range("D41").Value = "0.007191"
$ws.Range("E41").Value = "3.15%"
$ws.Range("G41").Value = "21"
$ws.Range("E42").Value = "-0.27%"
$ws.Range("G42").Value = "21"
$ws.Range("D43").Value = "0.002227"
$ws.Range("E43").Value = "1.56%"
$ws.Range("G43").Value = "21"
$ws.Range("D44").Value = "0.01076"
$ws.Range("E44").Value = "-15.68%"
$ws.Range("G44").Value = "21"
$ws.Range("D45").Value = "0.00006267"
$ws.Range("E45").Value = "1.69%"
$ws.Range("G45").Value = "21"
$ws.Range("E46").Value = "-33.33%"
$ws.Range("G46").Value = "21"
$ws.Range("D47").Value = "0.8255"
$ws.Range("E47").Value = "16.08%"
$ws.Range("G47").Value = "21"
$ws.Range("G48").Value = "21"
$ws.Range("G49").Value = "21"
$ws.Range("G50").Value = "21"
$ws.Range("G51").Value = "21"

# Restore the original (unformatted) style on the range so no stray
# number-format styling is left behind on the cells.
$dataRange.Style = "Normal"

